# Generate Report for Handback
#
# File "4821b8ad-595a-4a63-95d6-097abe830311.md" has finished its de-de
# handback (and was already done for zh-cn); update its status from
# "Ready for handoff" to "Handed back: in sync with en-US" on the
# Overview sheet (both the zh-cn and de-de columns) as well as on the
# per-language detail sheets, and stamp the "Latest Handback DateTime"
# for each language with the new handback timestamp.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the 4821b8ad...md file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $status   # zh-cn column
$overview.Range("C3").Value = $status   # de-de column

# --- zh-cn detail sheet: row 3 is the 4821b8ad...md file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status                  # Status
$zhcn.Range("H3").Value = "2016-03-21 16:49:50"     # Latest Handback DateTime

# --- de-de detail sheet: row 3 is the 4821b8ad...md file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status                  # Status
$dede.Range("H3").Value = "2016-03-21 16:49:56"     # Latest Handback DateTime
